# The workbook tracks weekly price observations for "Naranja" (orange)
# lots sold at "Terminal Hortofrutícola Agro Chillán". This commit adds
# a new weekly observation (two quality grades: Primera / Segunda for
# variety "Navel Late") as two new rows inserted right before the
# existing row 261, shifting all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 261 (old rows 261-302
# become 263-304).
$ws.Rows("261:262").Insert()

# --- New row 261: Navel Late / Primera, fecha 2021-11-22 (44522) ---
$ws.Cells.Item(261, 1).Value = 7
$ws.Cells.Item(261, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(261, 3).Value = "Ñuble"
$ws.Cells.Item(261, 4).Value = 44522
$ws.Cells.Item(261, 5).Value = 16
$ws.Cells.Item(261, 6).Value = "Fruta"
$ws.Cells.Item(261, 7).Value = 100102
$ws.Cells.Item(261, 8).Value = "Cítricos"
$ws.Cells.Item(261, 9).Value = 100102005
$ws.Cells.Item(261, 10).Value = "Naranja"
$ws.Cells.Item(261, 11).Value = "Navel Late"
$ws.Cells.Item(261, 12).Value = "Primera"
$ws.Cells.Item(261, 13).Value = 160
$ws.Cells.Item(261, 14).Value = 8000
$ws.Cells.Item(261, 15).Value = 9000
$ws.Cells.Item(261, 16).Value = 8500
$ws.Cells.Item(261, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(261, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(261, 19).Value = 567
$ws.Cells.Item(261, 20).Value = 15

# --- New row 262: Navel Late / Segunda, fecha 2021-11-22 (44522) ---
$ws.Cells.Item(262, 1).Value = 7
$ws.Cells.Item(262, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(262, 3).Value = "Ñuble"
$ws.Cells.Item(262, 4).Value = 44522
$ws.Cells.Item(262, 5).Value = 16
$ws.Cells.Item(262, 6).Value = "Fruta"
$ws.Cells.Item(262, 7).Value = 100102
$ws.Cells.Item(262, 8).Value = "Cítricos"
$ws.Cells.Item(262, 9).Value = 100102005
$ws.Cells.Item(262, 10).Value = "Naranja"
$ws.Cells.Item(262, 11).Value = "Navel Late"
$ws.Cells.Item(262, 12).Value = "Segunda"
$ws.Cells.Item(262, 13).Value = 50
$ws.Cells.Item(262, 14).Value = 7000
$ws.Cells.Item(262, 15).Value = 7000
$ws.Cells.Item(262, 16).Value = 7000
$ws.Cells.Item(262, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(262, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(262, 19).Value = 467
$ws.Cells.Item(262, 20).Value = 15
